$d = $word.ActiveDocument

# Remove the leading "test video logic" line (including its trailing manual line break)
# that precedes the "Seleccion..." sentence.
$find1 = "test video logic^l"
$d.Content.Find.Execute($find1, $false, $false, $false, $false, $false, $true, 1, $false, "", 2)

# Remove the trailing "acomodar barra hora-fecha" / "cambiar tiempos hora, clima" lines
# that follow "...Ticketera\Videos", replacing them with a single manual line break.
$find2 = "\Videos^lacomodar barra hora-fecha^lcambiar tiempos hora, clima"
$d.Content.Find.Execute($find2, $false, $false, $false, $false, $false, $true, 1, $false, "\Videos^l", 2)
